$wb = $excel.ActiveWorkbook

# ---- Sheet 1: LP1912 ----
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range("A2").Value = 'Última actualización: 07:57:27'
$ws1.Range("A3").Value = 'Total filas: 108'

# Swap C59/C60
$ws1.Cells.Item(59,3).Value = '16_SANTA ANA'
$ws1.Cells.Item(60,3).Value = '11_ETCHEVERRY'

# Rewrite rows 77-113 (new scrape data, chronologically merged)
$ws1.Cells.Item(77,1).Value = '07:57:27'
$ws1.Cells.Item(77,2).Value = '08:06'
$ws1.Cells.Item(77,3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(77,4).Value = 9
$ws1.Cells.Item(77,5).Value = 'LP1912'
$ws1.Cells.Item(78,1).Value = '07:57:27'
$ws1.Cells.Item(78,2).Value = '08:10'
$ws1.Cells.Item(78,3).Value = '16_SANTA ANA'
$ws1.Cells.Item(78,4).Value = 13
$ws1.Cells.Item(78,5).Value = 'LP1912'
$ws1.Cells.Item(79,1).Value = '07:44:08'
$ws1.Cells.Item(79,2).Value = '08:11'
$ws1.Cells.Item(79,3).Value = '16_SANTA ANA'
$ws1.Cells.Item(79,4).Value = 27
$ws1.Cells.Item(79,5).Value = 'LP1912'
$ws1.Cells.Item(80,1).Value = '06:37:24'
$ws1.Cells.Item(80,2).Value = '08:12'
$ws1.Cells.Item(80,3).Value = '15_ABASTO'
$ws1.Cells.Item(80,4).Value = 95
$ws1.Cells.Item(80,5).Value = 'LP1912'
$ws1.Cells.Item(81,1).Value = '07:44:08'
$ws1.Cells.Item(81,2).Value = '08:13'
$ws1.Cells.Item(81,3).Value = '10_OLMOS'
$ws1.Cells.Item(81,4).Value = 29
$ws1.Cells.Item(81,5).Value = 'LP1912'
$ws1.Cells.Item(82,1).Value = '06:37:24'
$ws1.Cells.Item(82,2).Value = '08:21'
$ws1.Cells.Item(82,3).Value = '26_HERNANDEZ'
$ws1.Cells.Item(82,4).Value = 104
$ws1.Cells.Item(82,5).Value = 'LP1912'
$ws1.Cells.Item(83,1).Value = '06:37:24'
$ws1.Cells.Item(83,2).Value = '08:22'
$ws1.Cells.Item(83,3).Value = '16_P MOR-SANTA ANA'
$ws1.Cells.Item(83,4).Value = 105
$ws1.Cells.Item(83,5).Value = 'LP1912'
$ws1.Cells.Item(84,1).Value = '06:37:24'
$ws1.Cells.Item(84,2).Value = '08:23'
$ws1.Cells.Item(84,3).Value = '215B_EL PATO'
$ws1.Cells.Item(84,4).Value = 106
$ws1.Cells.Item(84,5).Value = 'LP1912'
$ws1.Cells.Item(85,1).Value = '07:44:08'
$ws1.Cells.Item(85,2).Value = '08:23'
$ws1.Cells.Item(85,3).Value = '16_P MOR-SANTA ANA'
$ws1.Cells.Item(85,4).Value = 39
$ws1.Cells.Item(85,5).Value = 'LP1912'
$ws1.Cells.Item(86,1).Value = '06:37:24'
$ws1.Cells.Item(86,2).Value = '08:27'
$ws1.Cells.Item(86,3).Value = '84_COLONIA URQUIZA-ESC 49'
$ws1.Cells.Item(86,4).Value = 110
$ws1.Cells.Item(86,5).Value = 'LP1912'
$ws1.Cells.Item(87,1).Value = '07:57:27'
$ws1.Cells.Item(87,2).Value = '08:33'
$ws1.Cells.Item(87,3).Value = '10_OLMOS'
$ws1.Cells.Item(87,4).Value = 36
$ws1.Cells.Item(87,5).Value = 'LP1912'
$ws1.Cells.Item(88,1).Value = '06:52:38'
$ws1.Cells.Item(88,2).Value = '08:42'
$ws1.Cells.Item(88,3).Value = '81_EL PELIGRO'
$ws1.Cells.Item(88,4).Value = 110
$ws1.Cells.Item(88,5).Value = 'LP1912'
$ws1.Cells.Item(89,1).Value = '07:14:27'
$ws1.Cells.Item(89,2).Value = '08:43'
$ws1.Cells.Item(89,3).Value = '14_ABASTO'
$ws1.Cells.Item(89,4).Value = 89
$ws1.Cells.Item(89,5).Value = 'LP1912'
$ws1.Cells.Item(90,1).Value = '07:44:08'
$ws1.Cells.Item(90,2).Value = '08:44'
$ws1.Cells.Item(90,3).Value = '14_ABASTO'
$ws1.Cells.Item(90,4).Value = 60
$ws1.Cells.Item(90,5).Value = 'LP1912'
$ws1.Cells.Item(91,1).Value = '07:14:27'
$ws1.Cells.Item(91,2).Value = '08:54'
$ws1.Cells.Item(91,3).Value = '17_ROMERO'
$ws1.Cells.Item(91,4).Value = 100
$ws1.Cells.Item(91,5).Value = 'LP1912'
$ws1.Cells.Item(92,1).Value = '07:14:27'
$ws1.Cells.Item(92,2).Value = '09:01'
$ws1.Cells.Item(92,3).Value = '215A_EL PATO'
$ws1.Cells.Item(92,4).Value = 107
$ws1.Cells.Item(92,5).Value = 'LP1912'
$ws1.Cells.Item(93,1).Value = '07:44:08'
$ws1.Cells.Item(93,2).Value = '09:02'
$ws1.Cells.Item(93,3).Value = '215A_EL PATO'
$ws1.Cells.Item(93,4).Value = 78
$ws1.Cells.Item(93,5).Value = 'LP1912'
$ws1.Cells.Item(94,1).Value = '07:57:27'
$ws1.Cells.Item(94,2).Value = '09:03'
$ws1.Cells.Item(94,3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(94,4).Value = 66
$ws1.Cells.Item(94,5).Value = 'LP1912'
$ws1.Cells.Item(95,1).Value = '07:14:27'
$ws1.Cells.Item(95,2).Value = '09:07'
$ws1.Cells.Item(95,3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(95,4).Value = 113
$ws1.Cells.Item(95,5).Value = 'LP1912'
$ws1.Cells.Item(96,1).Value = '07:44:08'
$ws1.Cells.Item(96,2).Value = '09:08'
$ws1.Cells.Item(96,3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(96,4).Value = 84
$ws1.Cells.Item(96,5).Value = 'LP1912'
$ws1.Cells.Item(97,1).Value = '07:57:27'
$ws1.Cells.Item(97,2).Value = '09:09'
$ws1.Cells.Item(97,3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(97,4).Value = 72
$ws1.Cells.Item(97,5).Value = 'LP1912'
$ws1.Cells.Item(98,1).Value = '07:14:27'
$ws1.Cells.Item(98,2).Value = '09:10'
$ws1.Cells.Item(98,3).Value = '16_P MOR-SANTA ANA'
$ws1.Cells.Item(98,4).Value = 116
$ws1.Cells.Item(98,5).Value = 'LP1912'
$ws1.Cells.Item(99,1).Value = '07:44:08'
$ws1.Cells.Item(99,2).Value = '09:11'
$ws1.Cells.Item(99,3).Value = '16_P MOR-SANTA ANA'
$ws1.Cells.Item(99,4).Value = 87
$ws1.Cells.Item(99,5).Value = 'LP1912'
$ws1.Cells.Item(100,1).Value = '07:44:08'
$ws1.Cells.Item(100,2).Value = '09:14'
$ws1.Cells.Item(100,3).Value = '16_SANTA ANA'
$ws1.Cells.Item(100,4).Value = 90
$ws1.Cells.Item(100,5).Value = 'LP1912'
$ws1.Cells.Item(101,1).Value = '07:57:27'
$ws1.Cells.Item(101,2).Value = '09:16'
$ws1.Cells.Item(101,3).Value = '27_EL RETIRO'
$ws1.Cells.Item(101,4).Value = 79
$ws1.Cells.Item(101,5).Value = 'LP1912'
$ws1.Cells.Item(102,1).Value = '07:44:08'
$ws1.Cells.Item(102,2).Value = '09:17'
$ws1.Cells.Item(102,3).Value = '27_EL RETIRO'
$ws1.Cells.Item(102,4).Value = 93
$ws1.Cells.Item(102,5).Value = 'LP1912'
$ws1.Cells.Item(103,1).Value = '07:44:08'
$ws1.Cells.Item(103,2).Value = '09:21'
$ws1.Cells.Item(103,3).Value = '26_HERNANDEZ'
$ws1.Cells.Item(103,4).Value = 97
$ws1.Cells.Item(103,5).Value = 'LP1912'
$ws1.Cells.Item(104,1).Value = '07:57:27'
$ws1.Cells.Item(104,2).Value = '09:22'
$ws1.Cells.Item(104,3).Value = '16_SANTA ANA'
$ws1.Cells.Item(104,4).Value = 85
$ws1.Cells.Item(104,5).Value = 'LP1912'
$ws1.Cells.Item(105,1).Value = '07:57:27'
$ws1.Cells.Item(105,2).Value = '09:22'
$ws1.Cells.Item(105,3).Value = '17_ROMERO'
$ws1.Cells.Item(105,4).Value = 85
$ws1.Cells.Item(105,5).Value = 'LP1912'
$ws1.Cells.Item(106,1).Value = '07:57:27'
$ws1.Cells.Item(106,2).Value = '09:23'
$ws1.Cells.Item(106,3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(106,4).Value = 86
$ws1.Cells.Item(106,5).Value = 'LP1912'
$ws1.Cells.Item(107,1).Value = '07:44:08'
$ws1.Cells.Item(107,2).Value = '09:23'
$ws1.Cells.Item(107,3).Value = '17_ROMERO'
$ws1.Cells.Item(107,4).Value = 99
$ws1.Cells.Item(107,5).Value = 'LP1912'
$ws1.Cells.Item(108,1).Value = '07:44:08'
$ws1.Cells.Item(108,2).Value = '09:24'
$ws1.Cells.Item(108,3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(108,4).Value = 100
$ws1.Cells.Item(108,5).Value = 'LP1912'
$ws1.Cells.Item(109,1).Value = '07:44:08'
$ws1.Cells.Item(109,2).Value = '09:32'
$ws1.Cells.Item(109,3).Value = '15_ABASTO'
$ws1.Cells.Item(109,4).Value = 108
$ws1.Cells.Item(109,5).Value = 'LP1912'
$ws1.Cells.Item(110,1).Value = '07:44:08'
$ws1.Cells.Item(110,2).Value = '09:33'
$ws1.Cells.Item(110,3).Value = '10_OLMOS'
$ws1.Cells.Item(110,4).Value = 109
$ws1.Cells.Item(110,5).Value = 'LP1912'
$ws1.Cells.Item(111,1).Value = '07:44:08'
$ws1.Cells.Item(111,2).Value = '09:36'
$ws1.Cells.Item(111,3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(111,4).Value = 112
$ws1.Cells.Item(111,5).Value = 'LP1912'
$ws1.Cells.Item(112,1).Value = '07:44:08'
$ws1.Cells.Item(112,2).Value = '09:42'
$ws1.Cells.Item(112,3).Value = '215C_EL PATO'
$ws1.Cells.Item(112,4).Value = 118
$ws1.Cells.Item(112,5).Value = 'LP1912'
$ws1.Cells.Item(113,1).Value = '07:57:27'
$ws1.Cells.Item(113,2).Value = '09:43'
$ws1.Cells.Item(113,3).Value = '14_ABASTO'
$ws1.Cells.Item(113,4).Value = 106
$ws1.Cells.Item(113,5).Value = 'LP1912'

# ---- Sheet 2: LP1912-215 ----
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = 'Última actualización: 07:57:27'

# ---- Sheet 3: 6203-6173 ----
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = 'Última actualización: 07:57:27'
$ws3.Range("A3").Value = 'Total filas: 18'

# Insert new row at position 22 (push old row22 -> row23), then set rows 19-23
$ws3.Rows.Item(22).Insert()
$ws3.Cells.Item(19,1).Value = '07:57:27'
$ws3.Cells.Item(19,2).Value = '08:18'
$ws3.Cells.Item(19,3).Value = '215C_LA PLATA'
$ws3.Cells.Item(19,4).Value = 21
$ws3.Cells.Item(19,5).Value = 'L6203'
$ws3.Cells.Item(20,1).Value = '06:37:24'
$ws3.Cells.Item(20,2).Value = '08:30'
$ws3.Cells.Item(20,3).Value = '215A_LA PLATA'
$ws3.Cells.Item(20,4).Value = 113
$ws3.Cells.Item(20,5).Value = 'L6173'
$ws3.Cells.Item(21,1).Value = '06:52:38'
$ws3.Cells.Item(21,2).Value = '08:34'
$ws3.Cells.Item(21,3).Value = '215A_LA PLATA'
$ws3.Cells.Item(21,4).Value = 102
$ws3.Cells.Item(21,5).Value = 'L6173'
$ws3.Cells.Item(22,1).Value = '07:14:27'
$ws3.Cells.Item(22,2).Value = '08:35'
$ws3.Cells.Item(22,3).Value = '215A_LA PLATA'
$ws3.Cells.Item(22,4).Value = 81
$ws3.Cells.Item(22,5).Value = 'L6173'
$ws3.Cells.Item(23,1).Value = '07:14:27'
$ws3.Cells.Item(23,2).Value = '09:09'
$ws3.Cells.Item(23,3).Value = '215D_LA PLATA'
$ws3.Cells.Item(23,4).Value = 115
$ws3.Cells.Item(23,5).Value = 'L6203'
